# Edit: usuario.xlsx - rename/replace "None" password placeholders with real
# hash_password values (SHA-512 hex digests) for the first 8 users and a
# shared placeholder hash for the remaining users; change admin login to
# an email address and hyperlink it; adjust row heights for the now-taller
# wrapped SQL preview text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (password / hash_password) values -----------------------
$ws.Range("C3").Value  = "89b3bd6083da4e8c1b831b66a68a896cf77323a1a615afe859083d2c9ef9b6193e9d48ba4ac55eba00a60b47d967465c96fcc49b5ada64cca3c3fb36b35ac53b"
$ws.Range("C4").Value  = "274ad1a24014ff7f5102ace0fb916e479dd8900012ccecaf2279ae89b62c2bfbd0cf4b63c2697dbf6cec49cfb2dbcd8d95f0b1021ce70834a3f90bfa467a56f5"
$ws.Range("C5").Value  = "68dbb47e38b2c86c14ceebea4341b4fbc5aa0ec711a72fdbe7036be1ba93bfd8d2f453e53ef7411a8fabb4f6e5282c3c667d750e0e0e2f47162d7bb6d03a2261"
$ws.Range("C6").Value  = "6878baad9d3c064fa35754798a445810383a8914080097f9bf57fc69b736b86f19ece34f7f48af0425e6559c3daaec9654c7d8a4edcf80c8e90a336bc61e979a"
$ws.Range("C7").Value  = "7b0b1427b8a97db8f5c3fa7d2f597e5e836a9f6b9f552f6653ec4c8d29451a7e12af3fe79ae00c0c8a96e2833426b827b7cdaa2a525b59aff8b8f9623c519cbd"
$ws.Range("C8").Value  = "0b8ad5bf39b2db2a4d54625642d3ef3bff6760794d7d9c5a641888db830479938540f9cc958171af234faac0b67aeb500083ea1d7a6d8c96d107a5b6749fa190"

$restHash = "2a911471076d524988bf8512f67c215bbec5a40de9dd4ef2c2fd5c3d6cca2b4e12408b796498c28f052b922599b9afe6aa499062b00cf620d3bcab5ac3bbd850"
$ws.Range("C9").Value  = $restHash
$ws.Range("C10").Value = $restHash
$ws.Range("C11").Value = $restHash
$ws.Range("C12").Value = $restHash
$ws.Range("C13").Value = $restHash
$ws.Range("C14").Value = $restHash
$ws.Range("C15").Value = $restHash
$ws.Range("C16").Value = $restHash
$ws.Range("C17").Value = $restHash
$ws.Range("C18").Value = $restHash
$ws.Range("C19").Value = $restHash
$ws.Range("C20").Value = $restHash
$ws.Range("C21").Value = $restHash
$ws.Range("C22").Value = $restHash
$ws.Range("C23").Value = $restHash
$ws.Range("C24").Value = $restHash

# --- B3: admin login becomes an email address, with a mailto hyperlink -
$ws.Range("B3").Value = "admin@duoc.cl"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:admin@duoc.cl")

# --- Row heights: the longer wrapped SQL preview text grows row 3 and
# shrinks rows 4-24 slightly (font/wrap change) -------------------------
$ws.Rows.Item(3).RowHeight = 231
for ($r = 4; $r -le 24; $r++) {
    $ws.Rows.Item($r).RowHeight = 220.5
}

# --- View state: scroll/selection moved back to the top of the sheet ---
$ws.Range("B4").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
